$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.654.15"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.93"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.03"
$ws.Range("E5").Value = "  -1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.47"
$ws.Range("E6").Value = "  -5.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -7.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.50"
$ws.Range("E10").Value = "  -8.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.30"
$ws.Range("E12").Value = "  -7.18%  "
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.565.86"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.227.97"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  -4.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.91"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.569.28"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.04"
$ws.Range("E19").Value = "  -7.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0958"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  -5.39%  "
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.58"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  -7.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("E25").Value = "  -8.42%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.21"
$ws.Range("E29").Value = "  -6.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.94"
$ws.Range("E30").Value = "  -9.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.11"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.80"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0827"
$ws.Range("E33").Value = "  -5.79%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.22"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.88"
$ws.Range("E36").Value = "  -8.38%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.52"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.51"
$ws.Range("E40").Value = "  -10.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -11.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0306"
$ws.Range("E42").Value = "  -6.25%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.707.97"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.71"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.192"
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.10"
$ws.Range("E47").Value = "  -4.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.00"
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.05"
$ws.Range("E49").Value = "  -4.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.61"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.86"
$ws.Range("E51").Value = "  -6.34%  "

$ws.Range("D2,D3,D5,D6,D7,D9,D10,D11,D12,D14,D15,D16,D17,D18,D19,D20,D21,D23,D24,D25,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D39,D40,D41,D42,D44,D45,D46,D47,D48,D49,D50,D51").Style = "Normal"
